# Added Avoidance Mechanic Penalty.
# Adds a new "penalty" column (C) to the avoidance list sheet, with a
# penalty value for each existing mechanic_id / avoid_mechanic_id pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "penalty"

# Penalty values for rows 2-37 (aligned with existing data rows)
$penalties = @(
    22, 26, 27, 20, 30, 13, 28, 25, 30, 44,
    24, 14, 22, 42, 49, 31, 27, 46, 38, 20,
    41, 24, 29, 45, 50, 47, 21, 12, 20, 33,
    21, 13, 23, 21, 19, 42
)

for ($i = 0; $i -lt $penalties.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $penalties[$i]
}

# Reflect the post-edit selection on the new column.
$ws.Range("C2").Select()
